$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P11").Value = "Kenya"
$ws.Range("P12").Value = "Ethiopia"
$ws.Range("P13").Value = "Kenya"
$ws.Range("P14").Value = "United States"
$ws.Range("P15").Value = "Ethiopia"
